$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5590
$ws.Range("L3").Value = 6079
$ws.Range("L4").Value = 1499
$ws.Range("L6").Value = 4984
$ws.Range("L7").Value = 18514

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 202

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 360
$ws.Range("L3").Value = 434
$ws.Range("L6").Value = 304
$ws.Range("L7").Value = 1228

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 239
$ws.Range("L7").Value = 847

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 262

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 192
$ws.Range("L7").Value = 710

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 145
$ws.Range("L7").Value = 598
$ws.Range("L8").Value = 1228
$ws.Range("L11").Value = 304
$ws.Range("L18").Value = 126
$ws.Range("L19").Value = 506
$ws.Range("L20").Value = 456
$ws.Range("L23").Value = 204
$ws.Range("L25").Value = 114
$ws.Range("L29").Value = 1044
$ws.Range("L33").Value = 847
$ws.Range("L37").Value = 710
$ws.Range("L42").Value = 600
$ws.Range("L47").Value = 121
$ws.Range("L48").Value = 241
$ws.Range("L49").Value = 93
$ws.Range("L53").Value = 202
$ws.Range("L54").Value = 406
$ws.Range("L67").Value = 640
$ws.Range("L70").Value = 32
$ws.Range("L73").Value = 150
$ws.Range("L76").Value = 286
$ws.Range("L84").Value = 179
$ws.Range("L85").Value = 918
$ws.Range("L89").Value = 261
$ws.Range("L91").Value = 245
$ws.Range("L92").Value = 57
$ws.Range("L94").Value = 228
$ws.Range("L95").Value = 262
$ws.Range("L101").Value = 18514

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 247
$ws.Range("L4").Value = 42
$ws.Range("L7").Value = 640

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 58
$ws.Range("L7").Value = 179

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 73
$ws.Range("L4").Value = 34
$ws.Range("L6").Value = 194
$ws.Range("L7").Value = 406

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 401
$ws.Range("L4").Value = 56
$ws.Range("L6").Value = 263
$ws.Range("L7").Value = 1044

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 241

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 178
$ws.Range("L7").Value = 506

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 59
$ws.Range("L7").Value = 286

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 163
$ws.Range("L3").Value = 208
$ws.Range("L7").Value = 600

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 81
$ws.Range("L7").Value = 204

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 83
$ws.Range("L7").Value = 245

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 155
$ws.Range("L6").Value = 114
$ws.Range("L7").Value = 456

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 144
$ws.Range("L7").Value = 598

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 228

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L4").Value = 9
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 91
$ws.Range("L7").Value = 304

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 80
$ws.Range("L7").Value = 261

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 379
$ws.Range("L7").Value = 918
